$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "us-core-heart-rate"
$ws.Range("B2").Value = "US Core Heart Rate Profile"
$ws.Range("C2").Value = "null#vital-signs"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "LOINC#8867-4"
$ws.Range("F2").Value = "'"
$ws.Range("G2").Value = "dateTimeĵ, Periodĵ"
$ws.Range("H2").Value = "Quantityĵ"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "'"

$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
